# Lab1_ExpectedResultsTable.xlsx
# "Update Lab 1 Expected Results Table and register_test.sv for better clarity"
#
# Row 4 of Sheet1 gets its B4:G4 contents shifted one column to the right
# (so a freshly-typed 0 can be inserted at B4), and the sheet's saved
# selection moves from F11 to D28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# The cell that used to hold the big number-stored-as-text value (D4, text
# number format) ends up at E4, and the cell that used to hold a plain
# number (C4, general number format) ends up at D4. Move just the *formats*
# first (before any value is overwritten) so each destination cell keeps
# the correct number format/style without creating new style entries.
$ws.Range("D4").Copy()
$ws.Range("E4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Now shift the values themselves one column to the right. Go right-to-left
# so each source cell is read before it gets overwritten.
$ws.Range("G4").Value2 = $ws.Range("F4").Value2
$ws.Range("F4").Value2 = $ws.Range("E4").Value2
$ws.Range("E4").Value2 = $ws.Range("D4").Text
$ws.Range("D4").Value2 = $ws.Range("C4").Value2
$ws.Range("C4").Value2 = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = 0

# Update the sheet's stored selection.
$ws.Range("D28").Select()
